# fix: project management presentation
# Applies small copy-edit / grammar fixes across the deck and adds a
# missing bullet line to the "User Story" slide.

$p = $ppt.ActivePresentation

# --- Slide 2 : "Notion" text box -----------------------------------------
$s2 = $p.Slides.Item(2)
$sh = $s2.Shapes.Item(1)
$tr = $sh.TextFrame.TextRange
$tr.Paragraphs(1).Runs(1).Text = "Notion est une application de prise de notes collaboratives avec énormément de fonctionnalités"

# --- Slide 3 : Kanban overview captions ----------------------------------
$s3 = $p.Slides.Item(3)

$sh = $s3.Shapes.Item(3)
$tr = $sh.TextFrame.TextRange
$tr.Paragraphs(1).Runs(1).Text = "Les différentes stories peuvent être déplacés entre les différentes colonnes pour changer leur statut"

$sh = $s3.Shapes.Item(4)
$tr = $sh.TextFrame.TextRange
$tr.Paragraphs(1).Runs(1).Text = "Cela permet de voir facilement où en sont les différentes tâches et si quelqu’un est bloqué dessus"

# --- Slide 5 : User Story details ----------------------------------------
$s5 = $p.Slides.Item(5)
$sh = $s5.Shapes.Item(3)
$tr = $sh.TextFrame.TextRange
$tr.Paragraphs(1).Runs(1).Text = "Chaque User Story est une page à part entière qui détaille :"
$tr.Paragraphs(6).Runs(1).Text = "Ainsi que dans certains cas :"
$tr.Paragraphs(8).InsertAfter("`r- Les détails techniques")

# --- Slide 7 : Sprints -----------------------------------------------------
$s7 = $p.Slides.Item(7)
$sh = $s7.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange
$tr.Paragraphs(2).Runs(1).Text = "Terminés par un Sprint "
$tr.Paragraphs(4).Runs(1).Text = "Il y aura également des petites réunions journalières pour communiquer les avancées des différents membres de l’équipe "
